$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 35985.715
$ws.Range("I12").Value = 295.83334
$ws.Range("K12").Value = 295.83334
$ws.Range("M12").Value = -125.83334

$ws.Range("H33").Value = 1616.1666
$ws.Range("I33").Value = 239.4
$ws.Range("K33").Value = 239.4
$ws.Range("M33").Value = -10.40000000000001

$ws.Range("H43").Value = 1762.5834
$ws.Range("I43").Value = 3800.3333
$ws.Range("J43").Value = 1083.3334
$ws.Range("K43").Value = 3800.3333
$ws.Range("L43").Value = 1083.3334
$ws.Range("M43").Value = -3731.3333
$ws.Range("N43").Value = -1221.3334

$ws.Range("H94").Value = 3599.8
$ws.Range("I94").Value = 3599.8
$ws.Range("K94").Value = 3599.8
$ws.Range("M94").Value = -3148.8

$ws.Range("H132").Value = 3910140.2
$ws.Range("I132").Value = 4314394.5
$ws.Range("K132").Value = 12943183.5
$ws.Range("M132").Value = -12940653.5

$ws.Range("H137").Value = 1500.3846
$ws.Range("I137").Value = 1610.4736
$ws.Range("J137").Value = 1201.5714
$ws.Range("K137").Value = 4831.4208
$ws.Range("L137").Value = 3604.7142
$ws.Range("M137").Value = -2281.4208
$ws.Range("N137").Value = -8704.7142

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 2760
$ws.Range("I63").Value = 1520
$ws.Range("J63").Value = 4000
$ws.Range("K63").Value = 1520
$ws.Range("L63").Value = 4000
$ws.Range("M63").Value = -834
$ws.Range("N63").Value = -5372

$ws.Range("H66").Value = 2760
$ws.Range("I66").Value = 1520
$ws.Range("J66").Value = 4000
$ws.Range("K66").Value = 7600
$ws.Range("L66").Value = 20000
$ws.Range("M66").Value = -4168
$ws.Range("N66").Value = -26864

$ws.Range("H105").Value = 47000
$ws.Range("J105").Value = 47000
$ws.Range("L105").Value = 47000
$ws.Range("N105").Value = -53988

$ws.Range("H132").Value = 1787.9736
$ws.Range("J132").Value = 1803
$ws.Range("L132").Value = 5409
$ws.Range("N132").Value = -10469

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H62").Value = 48627
$ws.Range("J62").Value = 48627
$ws.Range("L62").Value = 48627
$ws.Range("N62").Value = -49999

$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()

$ws.Range("H65").Value = 48627
$ws.Range("J65").Value = 48627
$ws.Range("L65").Value = 145881
$ws.Range("N65").Value = -152745

$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 2671.2856
$ws.Range("I62").Value = 2599.5
$ws.Range("K62").Value = 2599.5
$ws.Range("M62").Value = -1975.5

$ws.Range("H65").Value = 2671.2856
$ws.Range("I65").Value = 2599.5
$ws.Range("K65").Value = 12997.5
$ws.Range("M65").Value = -9877.5

$ws.Range("H107").Value = 525.48486
$ws.Range("I107").Value = 489.16
$ws.Range("K107").Value = 489.16
$ws.Range("M107").Value = 1430.84

$ws.Range("H132").Value = 3132.6223
$ws.Range("I132").Value = 2911.5881
$ws.Range("K132").Value = 8734.764299999999
$ws.Range("M132").Value = -6204.764299999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H97").Value = 1231.3334
$ws.Range("I97").Value = 1099.6
$ws.Range("J97").Value = 1396
$ws.Range("K97").Value = 3298.8
$ws.Range("L97").Value = 4188
$ws.Range("M97").Value = -2802.8
$ws.Range("N97").Value = -5180

$ws.Range("H113").Value = 903.55554
$ws.Range("I113").Value = 1912.8572
$ws.Range("J113").Value = 550.3
$ws.Range("K113").Value = 5738.571599999999
$ws.Range("L113").Value = 1650.9
$ws.Range("M113").Value = -3568.571599999999
$ws.Range("N113").Value = -5990.9

$ws.Range("H138").Value = 2502
$ws.Range("I138").Value = 1743
$ws.Range("K138").Value = 5229
$ws.Range("M138").Value = -89

$ws.Range("H139").Value = 1874.4
$ws.Range("I139").Value = 902.0714
$ws.Range("J139").Value = 4143.1665
$ws.Range("K139").Value = 2706.2142
$ws.Range("L139").Value = 12429.4995
$ws.Range("M139").Value = 2433.7858
$ws.Range("N139").Value = -22709.4995

$ws.Range("H140").Value = 6481.4287
$ws.Range("J140").Value = 2408.3333
$ws.Range("L140").Value = 7224.999899999999
$ws.Range("N140").Value = -17584.9999

$ws.Range("H141").Value = 12346.667
$ws.Range("I141").Value = 12346.667
$ws.Range("K141").Value = 37040.001
$ws.Range("M141").Value = -31860.001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2969.5454
$ws.Range("I102").Value = 2032.8667
$ws.Range("J102").Value = 4976.7144
$ws.Range("K102").Value = 2032.8667
$ws.Range("L102").Value = 4976.7144
$ws.Range("M102").Value = -410.8667
$ws.Range("N102").Value = -8220.714400000001

$ws.Range("H122").Value = 1247.9166
$ws.Range("I122").Value = 955
$ws.Range("J122").Value = 1457.1428
$ws.Range("K122").Value = 2865
$ws.Range("L122").Value = 4371.428400000001
$ws.Range("M122").Value = -415
$ws.Range("N122").Value = -9271.428400000001

$ws.Range("H132").Value = 2709.0244
$ws.Range("I132").Value = 2585.3157
$ws.Range("K132").Value = 7755.9471
$ws.Range("M132").Value = -5225.9471

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1265762.4
$ws.Range("I46").Value = 445
$ws.Range("J46").Value = 1687534.9
$ws.Range("K46").Value = 445
$ws.Range("L46").Value = 1687534.9
$ws.Range("M46").Value = -257
$ws.Range("N46").Value = -1687910.9

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 9445.546
$ws.Range("J15").Value = 9490.1
$ws.Range("L15").Value = 9490.1
$ws.Range("N15").Value = -10066.1

$ws.Range("H107").Value = 333766.66
$ws.Range("I107").Value = 700
$ws.Range("J107").Value = 500300
$ws.Range("K107").Value = 2100
$ws.Range("L107").Value = 1500900
$ws.Range("M107").Value = -180
$ws.Range("N107").Value = -1504740

$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").ClearContents()

$ws.Range("H132").Value = 2566.5264
$ws.Range("I132").Value = 2786.8708
$ws.Range("J132").Value = 1590.7142
$ws.Range("K132").Value = 8360.6124
$ws.Range("L132").Value = 4772.142599999999
$ws.Range("M132").Value = -5830.6124
$ws.Range("N132").Value = -9832.142599999999
